# Applies the Lgi2-Adam23 NATMI update: drops the "Inflammatory-Mac" target-cluster rows
# (old sheet rows 17-19) and refreshes the TPM-derived NATMI metrics for every remaining
# Sending-cluster x Target-cluster combination (old rows 2-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 3 rows whose target cluster was "Inflammatory-Mac" (MuSCs sender) are gone entirely
# in the new data set, so remove them first (this also shrinks dimension to A1:T16).
$ws.Rows("17:19").Delete()

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.1376636666666667
$ws.Cells.Item(2,8).Value = 0.412991
$ws.Cells.Item(2,9).Value = 0.01821680097623009
$ws.Cells.Item(2,10).Value = 0.01821680097623009
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.2328916666666666
$ws.Cells.Item(2,14).Value = 0.6986749999999999
$ws.Cells.Item(2,15).Value = 0.01421300418632399
$ws.Cells.Item(2,16).Value = 0.01421300418632399
$ws.Cells.Item(2,17).Value = 0.03206072076944444
$ws.Cells.Item(2,18).Value = 0.288546486925
$ws.Cells.Item(2,19).Value = 0.0002589154685365892
$ws.Cells.Item(2,20).Value = 0.0002589154685365893

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.1376636666666667
$ws.Cells.Item(3,8).Value = 0.412991
$ws.Cells.Item(3,9).Value = 0.01821680097623009
$ws.Cells.Item(3,10).Value = 0.01821680097623009
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 14.57672866666667
$ws.Cells.Item(3,14).Value = 43.730186
$ws.Cells.Item(3,15).Value = 0.8895943273864486
$ws.Cells.Item(3,16).Value = 0.8895943273864487
$ws.Cells.Item(3,17).Value = 2.006685916258444
$ws.Cells.Item(3,18).Value = 18.060173246326
$ws.Cells.Item(3,19).Value = 0.0162055628115822
$ws.Cells.Item(3,20).Value = 0.01620556281158221

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.1376636666666667
$ws.Cells.Item(4,8).Value = 0.412991
$ws.Cells.Item(4,9).Value = 0.01821680097623009
$ws.Cells.Item(4,10).Value = 0.01821680097623009
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.423796666666667
$ws.Cells.Item(4,14).Value = 4.27139
$ws.Cells.Item(4,15).Value = 0.08689202268783405
$ws.Cells.Item(4,16).Value = 0.08689202268783405
$ws.Cells.Item(4,17).Value = 0.1960050697211111
$ws.Cells.Item(4,18).Value = 1.76404562749
$ws.Cells.Item(4,19).Value = 0.001582894683726342
$ws.Cells.Item(4,20).Value = 0.001582894683726342

# Row 5: ECs -> Neutrophils
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,4).Value = "Neutrophils"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.1376636666666667
$ws.Cells.Item(5,8).Value = 0.412991
$ws.Cells.Item(5,9).Value = 0.01821680097623009
$ws.Cells.Item(5,10).Value = 0.01821680097623009
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.136948
$ws.Cells.Item(5,14).Value = 0.410844
$ws.Cells.Item(5,15).Value = 0.008357716380185487
$ws.Cells.Item(5,16).Value = 0.008357716380185487
$ws.Cells.Item(5,17).Value = 0.01885276382266666
$ws.Cells.Item(5,18).Value = 0.169674874404
$ws.Cells.Item(5,19).Value = 0.0001522508559136172
$ws.Cells.Item(5,20).Value = 0.0001522508559136172

# Row 6: ECs -> Resolving-Mac
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.1376636666666667
$ws.Cells.Item(6,8).Value = 0.412991
$ws.Cells.Item(6,9).Value = 0.01821680097623009
$ws.Cells.Item(6,10).Value = 0.01821680097623009
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.01545066666666667
$ws.Cells.Item(6,14).Value = 0.046352
$ws.Cells.Item(6,15).Value = 0.0009429293592077716
$ws.Cells.Item(6,16).Value = 0.0009429293592077717
$ws.Cells.Item(6,17).Value = 0.002126995425777778
$ws.Cells.Item(6,18).Value = 0.019142958832
$ws.Cells.Item(6,19).Value = 0.00001717715647133215
$ws.Cells.Item(6,20).Value = 0.00001717715647133215

# Row 7: FAPs -> ECs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 7.314644999999999
$ws.Cells.Item(7,8).Value = 21.943935
$ws.Cells.Item(7,9).Value = 0.9679346439276632
$ws.Cells.Item(7,10).Value = 0.967934643927663
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.2328916666666666
$ws.Cells.Item(7,14).Value = 0.6986749999999999
$ws.Cells.Item(7,15).Value = 0.01421300418632399
$ws.Cells.Item(7,16).Value = 0.01421300418632399
$ws.Cells.Item(7,17).Value = 1.703519865124999
$ws.Cells.Item(7,18).Value = 15.331678786125
$ws.Cells.Item(7,19).Value = 0.0137572591462319
$ws.Cells.Item(7,20).Value = 0.0137572591462319

# Row 8: FAPs -> FAPs
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 7.314644999999999
$ws.Cells.Item(8,8).Value = 21.943935
$ws.Cells.Item(8,9).Value = 0.9679346439276632
$ws.Cells.Item(8,10).Value = 0.967934643927663
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 14.57672866666667
$ws.Cells.Item(8,14).Value = 43.730186
$ws.Cells.Item(8,15).Value = 0.8895943273864486
$ws.Cells.Item(8,16).Value = 0.8895943273864487
$ws.Cells.Item(8,17).Value = 106.62359545799
$ws.Cells.Item(8,18).Value = 959.6123591219099
$ws.Cells.Item(8,19).Value = 0.8610691685188712
$ws.Cells.Item(8,20).Value = 0.8610691685188712

# Row 9: FAPs -> MuSCs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 7.314644999999999
$ws.Cells.Item(9,8).Value = 21.943935
$ws.Cells.Item(9,9).Value = 0.9679346439276632
$ws.Cells.Item(9,10).Value = 0.967934643927663
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.423796666666667
$ws.Cells.Item(9,14).Value = 4.27139
$ws.Cells.Item(9,15).Value = 0.08689202268783405
$ws.Cells.Item(9,16).Value = 0.08689202268783405
$ws.Cells.Item(9,17).Value = 10.41456716885
$ws.Cells.Item(9,18).Value = 93.73110451964999
$ws.Cells.Item(9,19).Value = 0.08410579904050308
$ws.Cells.Item(9,20).Value = 0.08410579904050307

# Row 10: FAPs -> Neutrophils
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,4).Value = "Neutrophils"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 7.314644999999999
$ws.Cells.Item(10,8).Value = 21.943935
$ws.Cells.Item(10,9).Value = 0.9679346439276632
$ws.Cells.Item(10,10).Value = 0.967934643927663
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.136948
$ws.Cells.Item(10,14).Value = 0.410844
$ws.Cells.Item(10,15).Value = 0.008357716380185487
$ws.Cells.Item(10,16).Value = 0.008357716380185487
$ws.Cells.Item(10,17).Value = 1.00172600346
$ws.Cells.Item(10,18).Value = 9.015534031139998
$ws.Cells.Item(10,19).Value = 0.008089723228503237
$ws.Cells.Item(10,20).Value = 0.008089723228503235

# Row 11: FAPs -> Resolving-Mac
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 7.314644999999999
$ws.Cells.Item(11,8).Value = 21.943935
$ws.Cells.Item(11,9).Value = 0.9679346439276632
$ws.Cells.Item(11,10).Value = 0.967934643927663
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.01545066666666667
$ws.Cells.Item(11,14).Value = 0.046352
$ws.Cells.Item(11,15).Value = 0.0009429293592077716
$ws.Cells.Item(11,16).Value = 0.0009429293592077717
$ws.Cells.Item(11,17).Value = 0.11301614168
$ws.Cells.Item(11,18).Value = 1.01714527512
$ws.Cells.Item(11,19).Value = 0.000912693993553714
$ws.Cells.Item(11,20).Value = 0.000912693993553714

# Row 12: MuSCs -> ECs
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.104653
$ws.Cells.Item(12,8).Value = 0.313959
$ws.Cells.Item(12,9).Value = 0.01384855509610675
$ws.Cells.Item(12,10).Value = 0.01384855509610675
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.2328916666666666
$ws.Cells.Item(12,14).Value = 0.6986749999999999
$ws.Cells.Item(12,15).Value = 0.01421300418632399
$ws.Cells.Item(12,16).Value = 0.01421300418632399
$ws.Cells.Item(12,17).Value = 0.02437281159166666
$ws.Cells.Item(12,18).Value = 0.219355304325
$ws.Cells.Item(12,19).Value = 0.0001968295715555037
$ws.Cells.Item(12,20).Value = 0.0001968295715555037

# Row 13: MuSCs -> FAPs
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.104653
$ws.Cells.Item(13,8).Value = 0.313959
$ws.Cells.Item(13,9).Value = 0.01384855509610675
$ws.Cells.Item(13,10).Value = 0.01384855509610675
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 14.57672866666667
$ws.Cells.Item(13,14).Value = 43.730186
$ws.Cells.Item(13,15).Value = 0.8895943273864486
$ws.Cells.Item(13,16).Value = 0.8895943273864487
$ws.Cells.Item(13,17).Value = 1.525498385152667
$ws.Cells.Item(13,18).Value = 13.729485466374
$ws.Cells.Item(13,19).Value = 0.01231959605599526
$ws.Cells.Item(13,20).Value = 0.01231959605599526

# Row 14: MuSCs -> MuSCs
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,4).Value = "MuSCs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.104653
$ws.Cells.Item(14,8).Value = 0.313959
$ws.Cells.Item(14,9).Value = 0.01384855509610675
$ws.Cells.Item(14,10).Value = 0.01384855509610675
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 1.423796666666667
$ws.Cells.Item(14,14).Value = 4.27139
$ws.Cells.Item(14,15).Value = 0.08689202268783405
$ws.Cells.Item(14,16).Value = 0.08689202268783405
$ws.Cells.Item(14,17).Value = 0.1490045925566667
$ws.Cells.Item(14,18).Value = 1.34104133301
$ws.Cells.Item(14,19).Value = 0.001203328963604628
$ws.Cells.Item(14,20).Value = 0.001203328963604627

# Row 15: MuSCs -> Neutrophils
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,4).Value = "Neutrophils"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.104653
$ws.Cells.Item(15,8).Value = 0.313959
$ws.Cells.Item(15,9).Value = 0.01384855509610675
$ws.Cells.Item(15,10).Value = 0.01384855509610675
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.136948
$ws.Cells.Item(15,14).Value = 0.410844
$ws.Cells.Item(15,15).Value = 0.008357716380185487
$ws.Cells.Item(15,16).Value = 0.008357716380185487
$ws.Cells.Item(15,17).Value = 0.014332019044
$ws.Cells.Item(15,18).Value = 0.128988171396
$ws.Cells.Item(15,19).Value = 0.0001157422957686326
$ws.Cells.Item(15,20).Value = 0.0001157422957686326

# Row 16: MuSCs -> Resolving-Mac
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.104653
$ws.Cells.Item(16,8).Value = 0.313959
$ws.Cells.Item(16,9).Value = 0.01384855509610675
$ws.Cells.Item(16,10).Value = 0.01384855509610675
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.01545066666666667
$ws.Cells.Item(16,14).Value = 0.046352
$ws.Cells.Item(16,15).Value = 0.0009429293592077716
$ws.Cells.Item(16,16).Value = 0.0009429293592077717
$ws.Cells.Item(16,17).Value = 0.001616958618666667
$ws.Cells.Item(16,18).Value = 0.014552627568
$ws.Cells.Item(16,19).Value = 0.00001305820918272546
$ws.Cells.Item(16,20).Value = 0.00001305820918272546

